$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.5832667944600303
$ws.Range("D2").Value = 0.03024337164594471
$ws.Range("E2").Value = 0.2587809918563158
$ws.Range("F2").Value = 0.7407262362406257
$ws.Range("G2").Value = 0.002430470156977063
$ws.Range("I2").Value = 1.012226581589886
$ws.Range("K2").Value = 0.397417126174048
$ws.Range("L2").Value = 0.1535508338808711
$ws.Range("M2").Value = 0.1562473903505044
$ws.Range("N2").Value = 1.870339648719417
$ws.Range("O2").Value = 2.553903577800014

$ws.Range("B3").Value = 0.563823169258967
$ws.Range("D3").Value = 0.02829679280228703
$ws.Range("E3").Value = 0.2601573798439722
$ws.Range("F3").Value = 0.7376062893357584
$ws.Range("G3").Value = 0.002432801525739718
$ws.Range("I3").Value = 1.021250444267515
$ws.Range("K3").Value = 0.3473605290133435
$ws.Range("L3").Value = 0.1474363807469388
$ws.Range("M3").Value = 0.1506719846460562
$ws.Range("N3").Value = 1.887984375847335
$ws.Range("O3").Value = 2.556518631981845

$ws.Range("B4").Value = 0.5521071171004621
$ws.Range("D4").Value = 0.02709009284257746
$ws.Range("E4").Value = 0.2610546138234853
$ws.Range("F4").Value = 0.7361257566161115
$ws.Range("G4").Value = 0.002434310976002894
$ws.Range("I4").Value = 1.027190405284856
$ws.Range("K4").Value = 0.316495064848084
$ws.Range("L4").Value = 0.1437548336717569
$ws.Range("M4").Value = 0.1473138252807082
$ws.Range("N4").Value = 1.899384928679126
$ws.Range("O4").Value = 2.55953668177365

$ws.Range("B5").Value = 0.5473891563194968
$ws.Range("D5").Value = 0.02659548384782795
$ws.Range("E5").Value = 0.2614333814209526
$ws.Range("F5").Value = 0.7356319153387645
$ws.Range("G5").Value = 0.002434945755808738
$ws.Range("I5").Value = 1.029711458918303
$ws.Range("K5").Value = 0.3038849339450564
$ws.Range("L5").Value = 0.1422729634953015
$ws.Range("M5").Value = 0.1459618565250835
$ws.Range("N5").Value = 1.904173287347892
$ws.Range("O5").Value = 2.561121857602728

$ws.Range("B6").Value = 0.5466091651474443
$ws.Range("D6").Value = 0.02651318173839456
$ws.Range("E6").Value = 0.261497069842958
$ws.Range("F6").Value = 0.7355565282721628
$ws.Range("G6").Value = 0.002435052349764051
$ws.Range("I6").Value = 1.030136148258691
$ws.Range("K6").Value = 0.3017891043781731
$ws.Range("L6").Value = 0.1420280136271543
$ws.Range("M6").Value = 0.1457383638262293
$ws.Range("N6").Value = 1.904977001567683
$ws.Range("O6").Value = 2.561406537917279

$ws.Range("B7").Value = 0.5520432599057017
$ws.Range("D7").Value = 0.02708343396143675
$ws.Range("E7").Value = 0.2610596687822664
$ws.Range("F7").Value = 0.73611865309352
$ws.Range("G7").Value = 0.00243431945728308
$ws.Range("I7").Value = 1.027223998196927
$ws.Range("K7").Value = 0.3163251297903855
$ws.Range("L7").Value = 0.1437347740637165
$ws.Range("M7").Value = 0.1472955251521348
$ws.Range("N7").Value = 1.899448928935939
$ws.Range("O7").Value = 2.55955662124849

$ws.Range("B8").Value = 0.5765167625451113
$ws.Range("D8").Value = 0.02957458863693319
$ws.Range("E8").Value = 0.2592447696773601
$ws.Range("F8").Value = 0.7395602078047858
$ws.Range("G8").Value = 0.002431257862411231
$ws.Range("I8").Value = 1.015255183409355
$ws.Range("K8").Value = 0.3801850533446327
$ws.Range("L8").Value = 0.151427526290945
$ws.Range("M8").Value = 0.1543115377313065
$ws.Range("N8").Value = 1.876305986350268
$ws.Range("O8").Value = 2.554512181922973

$ws.Range("B9").Value = 0.6262549036084977
$ws.Range("D9").Value = 0.03436782872469735
$ws.Range("E9").Value = 0.2560979973662376
$ws.Range("F9").Value = 0.7497598846487179
$ws.Range("G9").Value = 0.002425870216236136
$ws.Range("I9").Value = 0.9949491647275792
$ws.Range("K9").Value = 0.5043573244998925
$ws.Range("L9").Value = 0.1670869183183328
$ws.Range("M9").Value = 0.168582364514176
$ws.Range("N9").Value = 1.835415293424427
$ws.Range("O9").Value = 2.555823807189682

$ws.Range("B10").Value = 0.6638402584531207
$ws.Range("D10").Value = 0.03783272945013039
$ws.Range("E10").Value = 0.2540354832552225
$ws.Range("F10").Value = 0.7593570197423105
$ws.Range("G10").Value = 0.002422283866647051
$ws.Range("I10").Value = 0.9819550761840929
$ws.Range("K10").Value = 0.5949221870722852
$ws.Range("L10").Value = 0.178938609814125
$ws.Range("M10").Value = 0.1793745752001072
$ws.Range("N10").Value = 1.808106192820444
$ws.Range("O10").Value = 2.563618123272818

$ws.Range("B11").Value = 0.6811609243762859
$ws.Range("D11").Value = 0.03939656713752981
$ws.Range("E11").Value = 0.2531509494560347
$ws.Range("F11").Value = 0.7641798147446792
$ws.Range("G11").Value = 0.002420732332621457
$ws.Range("I11").Value = 0.9764606976386396
$ws.Range("K11").Value = 0.6359743839250882
$ws.Range("L11").Value = 0.1844049493237492
$ws.Range("M11").Value = 0.1843500034236456
$ws.Range("N11").Value = 1.796274918834461
$ws.Range("O11").Value = 2.568647598264789

$ws.Range("B12").Value = 0.6877514092236652
$ws.Range("D12").Value = 0.0399869559635917
$ws.Range("E12").Value = 0.2528236931085202
$ws.Range("F12").Value = 0.7660717608687122
$ws.Range("G12").Value = 0.002420156240213037
$ws.Range("I12").Value = 0.9744399723649018
$ws.Range("K12").Value = 0.6514982356977441
$ws.Range("L12").Value = 0.1864856064785272
$ws.Range("M12").Value = 0.1862434493422782
$ws.Range("N12").Value = 1.791879764567607
$ws.Range("O12").Value = 2.570765447592976

$ws.Range("B13").Value = 0.6863306367392568
$ws.Range("D13").Value = 0.03985988573175803
$ws.Range("E13").Value = 0.2528938316963505
$ws.Range("F13").Value = 0.7656613772030596
$ws.Range("G13").Value = 0.002420279803877122
$ws.Range("I13").Value = 0.9748725093757784
$ws.Range("K13").Value = 0.6481558737927458
$ws.Range("L13").Value = 0.186037026708874
$ws.Range("M13").Value = 0.1858352477212222
$ws.Range("N13").Value = 1.792822555204102
$ws.Range("O13").Value = 2.57029984561288

$ws.Range("B14").Value = 0.6817024987576588
$ws.Range("D14").Value = 0.03944517505050271
$ws.Range("E14").Value = 0.2531238717637119
$ws.Range("F14").Value = 0.7643341507435082
$ws.Range("G14").Value = 0.002420684708283222
$ws.Range("I14").Value = 0.976293251587979
$ws.Range("K14").Value = 0.6372519799467966
$ws.Range("L14").Value = 0.1845759129343065
$ws.Range("M14").Value = 0.1845055915377714
$ws.Range("N14").Value = 1.795911621481711
$ws.Range("O14").Value = 2.568817560650729

$ws.Range("B15").Value = 0.6788717203637304
$ws.Range("D15").Value = 0.03919091729877522
$ws.Range("E15").Value = 0.2532657795895039
$ws.Range("F15").Value = 0.7635297350894916
$ws.Range("G15").Value = 0.002420934212102919
$ws.Range("I15").Value = 0.9771712936220673
$ws.Range("K15").Value = 0.6305701802234864
$ws.Range("L15").Value = 0.1836823254051438
$ws.Range("M15").Value = 0.1836923536857142
$ws.Range("N15").Value = 1.797814845246625
$ws.Range("O15").Value = 2.56793739359864

$ws.Range("B16").Value = 0.6627127523655929
$ws.Range("D16").Value = 0.03773027808854579
$ws.Range("E16").Value = 0.254094368205523
$ws.Range("F16").Value = 0.7590510299277469
$ws.Range("G16").Value = 0.002422386866683042
$ws.Range("I16").Value = 0.9823225279988073
$ws.Range("K16").Value = 0.5922363275560656
$ws.Range("L16").Value = 0.1785828708826642
$ws.Range("M16").Value = 0.1790507360976932
$ws.Range("N16").Value = 1.808891302000932
$ws.Range("O16").Value = 2.563319286977304

$ws.Range("B17").Value = 0.6528564630619371
$ws.Range("D17").Value = 0.03683103917806818
$ws.Range("E17").Value = 0.2546164189699283
$ws.Range("F17").Value = 0.756420498907211
$ws.Range("G17").Value = 0.002423298454125864
$ws.Range("I17").Value = 0.9855893383008691
$ws.Range("K17").Value = 0.5686818049051965
$ws.Range("L17").Value = 0.1754736421599716
$ws.Range("M17").Value = 0.1762200643438625
$ws.Range("N17").Value = 1.815837913178589
$ws.Range("O17").Value = 2.560866214319418

$ws.Range("B18").Value = 0.6472083980737011
$ws.Range("D18").Value = 0.03631265856125765
$ws.Range("E18").Value = 0.2549217461798303
$ws.Range("F18").Value = 0.7549505066049633
$ws.Range("G18").Value = 0.002423830300166995
$ws.Range("I18").Value = 0.9875075453337843
$ws.Range("K18").Value = 0.5551201460276047
$ws.Range("L18").Value = 0.1736923565183304
$ws.Range("M18").Value = 0.1745981566625012
$ws.Range("N18").Value = 1.819889120721818
$ws.Range("O18").Value = 2.559594938680988

$ws.Range("B19").Value = 0.6452996852244155
$ws.Range("D19").Value = 0.03613694493914466
$ws.Range("E19").Value = 0.2550259941574551
$ws.Range("F19").Value = 0.7544601829057598
$ws.Range("G19").Value = 0.002424011668046899
$ws.Range("I19").Value = 0.988163755406184
$ws.Range("K19").Value = 0.5505260661311127
$ws.Range("L19").Value = 0.1730904596158496
$ws.Range("M19").Value = 0.1740500788460793
$ws.Range("N19").Value = 1.821270356764387
$ws.Range("O19").Value = 2.559188496696322

$ws.Range("B20").Value = 0.6539035110544091
$ws.Range("D20").Value = 0.03692688515558729
$ws.Range("E20").Value = 0.254560322545083
$ws.Range("F20").Value = 0.7566960714859263
$ws.Range("G20").Value = 0.002423200635727147
$ws.Range("I20").Value = 0.9852375215341524
$ws.Range("K20").Value = 0.571190649125441
$ws.Range("L20").Value = 0.1758038948846661
$ws.Range("M20").Value = 0.1765207514033484
$ws.Range("N20").Value = 1.815092669798968
$ws.Range("O20").Value = 2.561112893611266

$ws.Range("B21").Value = 0.683061043712172
$ws.Range("D21").Value = 0.03956703474901957
$ws.Range("E21").Value = 0.253056094725711
$ws.Range("F21").Value = 0.7647222078566216
$ws.Range("G21").Value = 0.002420565468068593
$ws.Range("I21").Value = 0.9758743204751745
$ws.Range("K21").Value = 0.6404553130485624
$ws.Range("L21").Value = 0.1850047886367747
$ws.Range("M21").Value = 0.1848958908292886
$ws.Range("N21").Value = 1.795001978470129
$ws.Range("O21").Value = 2.569247155463046

$ws.Range("B22").Value = 0.7023006810103425
$ws.Range("D22").Value = 0.04128200681772398
$ws.Range("E22").Value = 0.2521178490364231
$ws.Range("F22").Value = 0.7703504305988105
$ws.Range("G22").Value = 0.002418909888811555
$ws.Range("I22").Value = 0.9701039139829533
$ws.Range("K22").Value = 0.6855968243671953
$ws.Range("L22").Value = 0.191080266503235
$ws.Range("M22").Value = 0.1904240203712178
$ws.Range("N22").Value = 1.782367534027962
$ws.Range("O22").Value = 2.575806584801455

$ws.Range("B23").Value = 0.6920155083671204
$ws.Range("D23").Value = 0.04036766458979457
$ws.Range("E23").Value = 0.2526145130098008
$ws.Range("F23").Value = 0.7673115469405616
$ws.Range("G23").Value = 0.00241978742062519
$ws.Range("I23").Value = 0.973151770673784
$ws.Range("K23").Value = 0.6615157974434567
$ws.Range("L23").Value = 0.1878320163031333
$ws.Range("M23").Value = 0.1874686109675707
$ws.Range("N23").Value = 1.789065392869514
$ws.Range("O23").Value = 2.572191963787276

$ws.Range("B24").Value = 0.6534300831569908
$ws.Range("D24").Value = 0.03688355753607198
$ws.Range("E24").Value = 0.2545856675610949
$ws.Range("F24").Value = 0.7565713532483116
$ws.Range("G24").Value = 0.002423244835245843
$ws.Range("I24").Value = 0.9853964530631032
$ws.Range("K24").Value = 0.5700564624322908
$ws.Range("L24").Value = 0.1756545681483033
$ws.Range("M24").Value = 0.1763847936944742
$ws.Range("N24").Value = 1.815429415194892
$ws.Range("O24").Value = 2.561000936832329

$ws.Range("B25").Value = 0.6126147513729734
$ws.Range("D25").Value = 0.03308103778465465
$ws.Range("E25").Value = 0.2569053459641697
$ws.Range("F25").Value = 0.7466311940704671
$ws.Range("G25").Value = 0.002427262135641344
$ws.Range("I25").Value = 1.000104144970759
$ws.Range("K25").Value = 0.4708806644285346
$ws.Range("L25").Value = 0.1627895187129553
$ws.Range("M25").Value = 0.1646673796856852
$ws.Range("N25").Value = 1.845996928474793
$ws.Range("O25").Value = 2.554269476468647
